$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = '34.383.06'
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value = '  +0.59%  '

$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = '1.791.42'
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value = '  +0.36%  '

$ws.Cells.Item(4,5).Value = '  -0.13%  '

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '226.26'
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = '  +0.02%  '

$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '0.553'
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value = '  +0.97%  '

$ws.Cells.Item(7,5).Value = '  -0.12%  '

$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = '32.67'
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value = '  +1.71%  '

$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '0.296'
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = '  +1.30%  '

$ws.Cells.Item(10,5).Value = '  +0.28%  '

$ws.Cells.Item(11,5).Value = '  -0.59%  '

$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = '2.048.49'
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value = '  +0.30%  '

$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '1.796.05'
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value = '  +1.23%  '

$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '11.07'
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = '  +0.98%  '

$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '0.634'
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value = '  +1.29%  '

$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = '34.377.90'
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value = '  +0.59%  '

$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '4.28'
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Value = '  +2.36%  '

$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = '68.36'
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value = '  +0.98%  '

$ws.Cells.Item(19,5).Value = '  +0.31%  '

$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '244.30'
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = '  -0.64%  '

$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '11.25'
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value = '  +2.50%  '

$ws.Cells.Item(22,5).Value = '  +0.01%  '

$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '4.16'
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value = '  +0.70%  '

$ws.Cells.Item(24,5).Value = '  +1.30%  '

$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = '166.47'
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value = '  +2.98%  '

$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '7.30'
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value = '  +2.11%  '

$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '16.50'
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = '  +1.14%  '

$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = '0.116'
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Value = '  +0.99%  '

$ws.Cells.Item(29,5).Value = '  -0.24%  '

$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = '3.98'
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value = '  +6.44%  '

$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = '0.0525'
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).Value = '  +1.18%  '

$ws.Cells.Item(32,5).Value = '  +1.71%  '

$ws.Cells.Item(33,5).Value = '  +0.56%  '

$ws.Cells.Item(34,5).Value = '  +1.09%  '

$ws.Cells.Item(35,5).Value = '  -0.51%  '

$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = '1.401.78'
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).Value = '  -3.08%  '

$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '0.676'
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Value = '  +3.09%  '

$ws.Cells.Item(38,5).Value = '  +2.12%  '

$ws.Cells.Item(39,5).Value = '  -0.50%  '

$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '85.04'
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = '  +2.32%  '

$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '2.83'
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = '  +4.29%  '

$ws.Cells.Item(42,5).Value = '  +1.05%  '

$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '0.937'
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value = '  +2.29%  '

$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '13.83'
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value = '  +0.44%  '

$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '0.0526'
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = '  +2.18%  '

$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '1.11'
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value = '  +3.05%  '

$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '6.02'
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value = '  -1.03%  '

$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = '1.950.12'
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value = '  +0.39%  '

$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '105.06'
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = '  +0.19%  '

$ws.Cells.Item(50,5).Value = '  -0.11%  '

$ws.Cells.Item(51,5).Value = '  -2.44%  '
